# Insert a new weekly price-report row at position 169, shifting the
# existing rows 169-271 down to 170-272 (dimension grows from R271 to R272).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value  = 10
$ws.Cells.Item(169, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value  = "La Araucanía"
$ws.Cells.Item(169, 4).Value  = 44606
$ws.Cells.Item(169, 5).Value  = 9
$ws.Cells.Item(169, 6).Value  = 100112044
$ws.Cells.Item(169, 7).Value  = "Perejil"
$ws.Cells.Item(169, 8).Value  = "Sin especificar"
$ws.Cells.Item(169, 9).Value  = "Primera"
$ws.Cells.Item(169, 10).Value = 80
$ws.Cells.Item(169, 11).Value = 4000
$ws.Cells.Item(169, 12).Value = 4000
$ws.Cells.Item(169, 13).Value = 4000
$ws.Cells.Item(169, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(169, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(169, 16).Value = 1333
$ws.Cells.Item(169, 17).Value = 3
$ws.Cells.Item(169, 18).Value = "Hortaliza"
